$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.802.87'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').Value = '3.487.90'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'583.24"
$ws.Range('E5').Value = '  -2.00%  '
$ws.Range('D6').Value = "'129.91"
$ws.Range('E6').Value = '  -3.17%  '
$ws.Range('D7').Value = '3.490.11'
$ws.Range('E7').Value = '  -0.68%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = "'0.482"
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').Value = "'7.10"
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').Value = "'0.377"
$ws.Range('D13').Value = '4.062.54'
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('D14').Value = "'27.18"
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').Value = '3.505.08'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = "'0.0000176"
$ws.Range('E17').Value = '  -3.22%  '
$ws.Range('D18').Value = '63.837.94'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('D19').Value = "'9.83"
$ws.Range('E19').Value = '  -1.85%  '
$ws.Range('D20').Value = "'14.05"
$ws.Range('E20').Value = '  -1.99%  '
$ws.Range('E21').Value = '  -1.51%  '
$ws.Range('D22').Value = "'379.91"
$ws.Range('E22').Value = '  -3.15%  '
$ws.Range('E23').Value = '  -1.46%  '
$ws.Range('D24').Value = '3.622.16'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').Value = "'73.18"
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').Value = "'0.0000113"
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('E28').Value = '  -2.16%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').Value = "'7.43"
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('D31').Value = "'8.19"
$ws.Range('E31').Value = '  -1.83%  '
$ws.Range('E32').Value = '  -2.91%  '
$ws.Range('D33').Value = '3.492.11'
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('D35').Value = "'23.32"
$ws.Range('E35').Value = '  -3.64%  '
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('D37').Value = "'5.26"
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('E39').Value = '  -1.92%  '
$ws.Range('D40').Value = "'160.21"
$ws.Range('E40').Value = '  -4.96%  '
$ws.Range('D41').Value = "'0.0790"
$ws.Range('E41').Value = '  -3.42%  '
$ws.Range('E42').Value = '  -1.75%  '
$ws.Range('D43').Value = "'26.13"
$ws.Range('E43').Value = '  +1.69%  '
$ws.Range('D44').Value = "'1.00"
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('E46').Value = '  -3.93%  '
$ws.Range('D47').Value = "'4.35"
$ws.Range('E47').Value = '  -1.84%  '
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('D49').Value = '2.419.57'
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('D50').Value = "'6.80"
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('D51').Value = "'0.887"
$ws.Range('E51').Value = '  -1.31%  '
